$wb = $excel.ActiveWorkbook

# --- status_sheet1: reorder/update the first four replies and their statuses,
#     and clear the previously-set statuses at rows 49-51 ---
$ws1 = $wb.Worksheets.Item("status_sheet1")
$ws1.Range("B2").Value = "you can contact at the end of this month"
$ws1.Range("C2").Value = "Rescheduled"
$ws1.Range("B3").Value = "anytime after this week"
$ws1.Range("C3").Value = "Rescheduled"
$ws1.Range("B4").Value = "yeah tell"
$ws1.Range("C4").Value = "Paragraph Repeat"
$ws1.Range("B5").Value = "ok bye"
$ws1.Range("C5").Value = "Passed"
$ws1.Range("C49").ClearContents()
$ws1.Range("C50").ClearContents()
$ws1.Range("C51").ClearContents()

# --- status_sheet2: set C5 status, clear C49 and C51 ---
$ws2 = $wb.Worksheets.Item("status_sheet2")
$ws2.Range("C5").Value = "Passed"
$ws2.Range("C49").ClearContents()
$ws2.Range("C51").ClearContents()

# --- status_sheet3: set C5 status, clear C51 ---
$ws3 = $wb.Worksheets.Item("status_sheet3")
$ws3.Range("C5").Value = "Passed"
$ws3.Range("C51").ClearContents()

# --- status_sheet4: set C5 status, clear C51 ---
$ws4 = $wb.Worksheets.Item("status_sheet4")
$ws4.Range("C5").Value = "Passed"
$ws4.Range("C51").ClearContents()

# --- status_sheet5: set C5 status, clear C51 ---
$ws5 = $wb.Worksheets.Item("status_sheet5")
$ws5.Range("C5").Value = "Passed"
$ws5.Range("C51").ClearContents()

# --- status_sheet6: set C5 status, clear C51 ---
$ws6 = $wb.Worksheets.Item("status_sheet6")
$ws6.Range("C5").Value = "Passed"
$ws6.Range("C51").ClearContents()
